# Update avg_long (U) / avg_short (V) swap values per row, per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: XAUUSD
$ws.Cells.Item(4, 21).Value = -3.2329
$ws.Cells.Item(4, 22).Value = -2.3793
# Row 7: .A50
$ws.Cells.Item(7, 21).Value = 31.13
$ws.Cells.Item(7, 22).Value = -57.728
# Row 8: .AUS200
$ws.Cells.Item(8, 21).Value = -3.2294
$ws.Cells.Item(8, 22).Value = -7.4361
# Row 9: .DE30
$ws.Cells.Item(9, 21).Value = -80.2813
$ws.Cells.Item(9, 22).Value = -160.2216
# Row 10: .ES35
$ws.Cells.Item(10, 21).Value = 19.7708
$ws.Cells.Item(10, 22).Value = -40.9614
# Row 11: .F40
$ws.Cells.Item(11, 21).Value = 4.4625
$ws.Cells.Item(11, 22).Value = -15.5193
# Row 12: .HK50
$ws.Cells.Item(12, 21).Value = 224.688
$ws.Cells.Item(12, 22).Value = -318.362
# Row 13: .JP225
$ws.Cells.Item(13, 21).Value = -16.9304
$ws.Cells.Item(13, 22).Value = -24.9555
# Row 14: .STOXX50
$ws.Cells.Item(14, 21).Value = 0.5618
$ws.Cells.Item(14, 22).Value = -6.9554
# Row 15: .UK100
$ws.Cells.Item(15, 21).Value = 1.5068
$ws.Cells.Item(15, 22).Value = -12.8775
# Row 17: .US100
$ws.Cells.Item(17, 21).Value = -7.603
$ws.Cells.Item(17, 22).Value = -13.4897
# Row 18: .US30
$ws.Cells.Item(18, 21).Value = -17.0972
$ws.Cells.Item(18, 22).Value = -37.4958
# Row 19: .US500
$ws.Cells.Item(19, 21).Value = -14.4655
$ws.Cells.Item(19, 22).Value = -48.4154
# Row 20: .USOil
$ws.Cells.Item(20, 21).Value = -2.01
$ws.Cells.Item(20, 22).Value = -22.935
# Row 21: .XNGUSD
$ws.Cells.Item(21, 21).Value = -1.87
$ws.Cells.Item(21, 22).Value = 0.3654
# Row 22: AUDCAD
$ws.Cells.Item(22, 21).Value = -2.1549
$ws.Cells.Item(22, 22).Value = -0.8118
# Row 23: AUDCHF
$ws.Cells.Item(23, 21).Value = 0.4536
$ws.Cells.Item(23, 22).Value = -2.9612
# Row 24: AUDJPY
$ws.Cells.Item(24, 21).Value = -0.9031
$ws.Cells.Item(24, 22).Value = -1.5807
# Row 25: AUDNZD
$ws.Cells.Item(25, 21).Value = -2.442
$ws.Cells.Item(25, 22).Value = -0.9735
# Row 26: AUDSGD
$ws.Cells.Item(26, 21).Value = -3.3198
$ws.Cells.Item(26, 22).Value = -1.2551
# Row 27: AUDUSD
$ws.Cells.Item(27, 21).Value = -2.6125
$ws.Cells.Item(27, 22).Value = -1.8117
# Row 28: CADCHF
$ws.Cells.Item(28, 21).Value = 0.837
$ws.Cells.Item(28, 22).Value = -3.6971
# Row 29: CADJPY
$ws.Cells.Item(29, 21).Value = -0.3311
$ws.Cells.Item(29, 22).Value = -2.299
# Row 30: CADSGD
$ws.Cells.Item(30, 21).Value = 0.0496
$ws.Cells.Item(30, 22).Value = -1.7589
# Row 31: CHFJPY
$ws.Cells.Item(31, 21).Value = -4.4539
$ws.Cells.Item(31, 22).Value = 0.6173999999999999
# Row 32: CHFSGD
$ws.Cells.Item(32, 21).Value = -9.9
$ws.Cells.Item(32, 22).Value = -2.7588
# Row 33: EURAUD
$ws.Cells.Item(33, 21).Value = -5.8729
$ws.Cells.Item(33, 22).Value = 0.2845
# Row 34: EURCAD
$ws.Cells.Item(34, 21).Value = -6.1072
$ws.Cells.Item(34, 22).Value = 0.9324
# Row 35: EURCHF
$ws.Cells.Item(35, 21).Value = -0.8756
$ws.Cells.Item(35, 22).Value = -2.574
# Row 37: EURGBP
$ws.Cells.Item(37, 21).Value = -4.0447
# Row 39: EURJPY
$ws.Cells.Item(39, 21).Value = -4.2273
$ws.Cells.Item(39, 22).Value = -0.0001
# Row 42: EURNZD
$ws.Cells.Item(42, 21).Value = -7.029
$ws.Cells.Item(42, 22).Value = 1.0467
# Row 45: EURSGD
$ws.Cells.Item(45, 21).Value = -8.1972
$ws.Cells.Item(45, 22).Value = 0.3348
# Row 47: EURUSD
$ws.Cells.Item(47, 21).Value = -5.3967
$ws.Cells.Item(47, 22).Value = 0.5833
# Row 48: GBPAUD
$ws.Cells.Item(48, 21).Value = -3.8797
$ws.Cells.Item(48, 22).Value = -4.8565
# Row 49: GBPCAD
$ws.Cells.Item(49, 21).Value = -4.9852
$ws.Cells.Item(49, 22).Value = -3.3242
# Row 50: GBPCHF
$ws.Cells.Item(50, 21).Value = 0.3789
$ws.Cells.Item(50, 22).Value = -6.8959
# Row 52: GBPJPY
$ws.Cells.Item(52, 21).Value = -2.4926
$ws.Cells.Item(52, 22).Value = -4.5177
# Row 54: GBPNZD
$ws.Cells.Item(54, 21).Value = -6.162
$ws.Cells.Item(54, 22).Value = -4.3116
# Row 56: GBPSGD
$ws.Cells.Item(56, 21).Value = -5.4868
$ws.Cells.Item(56, 22).Value = -2.6708
# Row 57: GBPUSD
$ws.Cells.Item(57, 21).Value = -3.8808
$ws.Cells.Item(57, 22).Value = -3.0316
# Row 62: NZDCAD
$ws.Cells.Item(62, 21).Value = -1.3464
$ws.Cells.Item(62, 22).Value = -1.3519
# Row 63: NZDCHF
$ws.Cells.Item(63, 21).Value = 0.6912
$ws.Cells.Item(63, 22).Value = -3.2351
# Row 64: NZDJPY
$ws.Cells.Item(64, 21).Value = -0.3113
$ws.Cells.Item(64, 22).Value = -1.9767
# Row 65: NZDSGD
$ws.Cells.Item(65, 21).Value = -2.4409
$ws.Cells.Item(65, 22).Value = -1.87
# Row 75: USDJPY
$ws.Cells.Item(75, 21).Value = -1.3686
$ws.Cells.Item(75, 22).Value = -3.5619
